$wb = $excel.ActiveWorkbook

# Each entry: target worksheet, row number, and the new values for columns H..N
# (a value of $null means that column is left unchanged for that row)
$updates = @(
  @{Sheet="ALC"; Row=76; Vals=@(27030332, 30306170, 4676, 30306170, 4676, -30305855, -5306)},
  @{Sheet="ALC"; Row=79; Vals=@(27030332, 30306170, 4676, 30306170, 4676, -30305078, -6860)},
  @{Sheet="ALC"; Row=100; Vals=@(6425.8335, 2615.7144, 11760, 2615.7144, 11760, -2074.7144, -12842)},
  @{Sheet="ALC"; Row=129; Vals=@(618.7241, 497.4091, $null, 1492.2273, $null, 3507.7727, $null)},
  @{Sheet="ALC"; Row=132; Vals=@(170411.6, 4089.1538, 479296.16, 12267.4614, 1437888.48, -9737.4614, -1442948.48)},
  @{Sheet="ALC"; Row=135; Vals=@(7693186, 253.97368, 18520276, 2285.76312, 166682484, 249.2368799999999, -166687554)},
  @{Sheet="ALC"; Row=137; Vals=@(14859.443, 18111.639, 5877.1904, 54334.917, 17631.5712, -51784.917, -22731.5712)},
  @{Sheet="ALC"; Row=138; Vals=@(1203.85, 611.4655, 2021.9048, 1834.3965, 6065.7144, 3305.6035, -16345.7144)},
  @{Sheet="ALC"; Row=141; Vals=@(1371.4902, 909.7857, 3526.111, 2729.3571, 10578.333, 2450.6429, -20938.333)},
  @{Sheet="ARM"; Row=32; Vals=@(4731.4546, 4694.944, 5056.4, 4694.944, 5056.4, -4407.944, -5630.4)},
  @{Sheet="ARM"; Row=61; Vals=@(1566.151, 1697.8, 1161.0769, 1697.8, 1161.0769, -1485.8, -1585.0769)},
  @{Sheet="ARM"; Row=74; Vals=@(18426.846, 24429.465, 1219.3334, 24429.465, 1219.3334, -23555.465, -2967.3334)},
  @{Sheet="ARM"; Row=77; Vals=@(18426.846, 24429.465, 1219.3334, 122147.325, 6096.666999999999, -117779.325, -14832.667)},
  @{Sheet="ARM"; Row=132; Vals=@(1275905.6, 1480612.6, 532495.75, 4441837.800000001, 1597487.25, -4439307.800000001, -1602547.25)},
  @{Sheet="ARM"; Row=136; Vals=@(1566.151, 1697.8, 1161.0769, 5093.4, 3483.2307, -2543.4, -8583.2307)},
  @{Sheet="BSM"; Row=86; Vals=@(390689.28, 1850, 876738.4, 1850, 876738.4, -727, -878984.4)},
  @{Sheet="BSM"; Row=89; Vals=@(390689.28, 1850, 876738.4, 9250, 4383692, -3634, -4394924)},
  @{Sheet="BSM"; Row=107; Vals=@(446.66666, 416, 600, 416, 600, 1504, -4440)},
  @{Sheet="BSM"; Row=134; Vals=@(19147.918, 905.55817, 62726.89, 2716.67451, 188180.67, -181.6745099999998, -193250.67)},
  @{Sheet="CRP"; Row=31; Vals=@(204947.17, 250809.78, 16910.5, 250809.78, 16910.5, -250514.78, -17500.5)},
  @{Sheet="CRP"; Row=34; Vals=@(204947.17, 250809.78, 16910.5, 250809.78, 16910.5, -250607.78, -17314.5)},
  @{Sheet="CRP"; Row=58; Vals=@(2133.7258, 740.06384, 6500.533, 740.06384, 6500.533, -537.06384, -6906.533)},
  @{Sheet="CRP"; Row=132; Vals=@(961.1964, 996.4286, 714.5714, 2989.2858, 2143.7142, -459.2857999999997, -7203.7142)},
  @{Sheet="CRP"; Row=136; Vals=@(2133.7258, 740.06384, 6500.533, 2220.19152, 19501.599, 329.8084799999997, -24601.599)},
  @{Sheet="CUL"; Row=131; Vals=@(20492630, 464.2857, 23149022, 1392.8571, 69447066, 3647.1429, -69457146)},
  @{Sheet="CUL"; Row=132; Vals=@(1233.4546, 1202, 1375, 10818, 12375, -8288, -17435)},
  @{Sheet="GSM"; Row=70; Vals=@(3850366.5, 5003879.5, 5323, 5003879.5, 5323, -5003609.5, -5863)},
  @{Sheet="GSM"; Row=73; Vals=@(3850366.5, 5003879.5, 5323, 5003879.5, 5323, -5002943.5, -7195)},
  @{Sheet="GSM"; Row=132; Vals=@(23760.844, 1337.1923, 54445.844, 4011.5769, 163337.532, -1481.5769, -168397.532)},
  @{Sheet="LTW"; Row=136; Vals=@(205461.12, 264161.62, 2677.5454, 792484.86, 8032.6362, -789934.86, -13132.6362)},
  @{Sheet="WVR"; Row=96; Vals=@(2323.647, 1775.375, 2811, 1775.375, 2811, -402.375, -5557)},
  @{Sheet="WVR"; Row=136; Vals=@(1358044.2, 1553275.2, 716571, 4659825.6, 2149713, -4657275.6, -2154813)}

)

$cols = @("H", "I", "J", "K", "L", "M", "N")

$sheetsTouched = @{}
foreach ($entry in $updates) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    $sheetsTouched[$entry.Sheet] = $true
    $rowVals = $entry.Vals
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $v = $rowVals[$i]
        if ($v -ne $null) {
            $ws.Range($cols[$i] + $entry.Row).Value = $v
        }
    }
}

Write-Host ("Updated " + $updates.Length + " rows across " + $sheetsTouched.Keys.Count + " sheets")
